# Contact List App - RTM: re-sync the "Bug ID" column (J) against the
# reshuffled bug numbering, fix the "Current Staus" header typo, flip
# row 14's test result to Passed (with matching Passed styling), and
# restore the row heights that shrank now that J6/J9 hold shorter text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header typo fix: "Current Staus" -> "Current Status"
$ws.Range("I5").Value = "Current Status"

# Re-grouped Bug ID values for each requirement row (column J)
$ws.Range("J6").Value  = "Bug - 001 / Bug - 002 / Bug - 003 / Bug - 004 / Bug - 005 / Bug - 006"
$ws.Range("J7").Value  = "Bug - 007"
$ws.Range("J8").Value  = "Bug - 008"
$ws.Range("J9").Value  = "Bug - 009 / Bug - 010 / Bug - 011 / Bug - 012 / Bug - 013 / Bug - 014"
$ws.Range("J10").Value = "Bug - 015"
$ws.Range("J11").Value = "Bug - 016 / Bug - 017 / Bug - 018 / Bug - 019 / Bug - 020"
$ws.Range("J12").Value = "Bug - 027 / Bug - 028"
$ws.Range("J13").Value = "Bug - 021 / Bug - 022 / Bug - 023"
$ws.Range("J14").Value = "-"
$ws.Range("J15").Value = "-"
$ws.Range("J16").Value = "Bug - 024"
$ws.Range("J17").Value = "Bug - 025 / Bug - 026 / Bug - 027 / Bug - 028"

# Row 14's test status flips from Failed to Passed - copy the "Passed"
# formatting (green fill) from I15 so the cell matches visually too.
$ws.Range("I15").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = "Passed"

# J6 / J9 now hold much shorter text, so their wrapped row height shrinks.
$ws.Rows("6:6").RowHeight = 43.5
$ws.Rows("9:9").RowHeight = 43.5

# Restore the scroll/selection state captured in the workbook view.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J15").Select()
